$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44351
$ws.Range("N2").Value = 15000
$ws.Range("O2").Value = 15000
$ws.Range("P2").Value = 15000
$ws.Range("S2").Value = 1000
$ws.Range("D3").Value = 44351
$ws.Range("M3").Value = 200
$ws.Range("D4").Value = 44309
$ws.Range("N4").Value = 17500
$ws.Range("O4").Value = 17500
$ws.Range("P4").Value = 17500
$ws.Range("S4").Value = 1167
$ws.Range("D5").Value = 44309
$ws.Range("D6").Value = 44344
$ws.Range("M6").Value = 100
$ws.Range("N6").Value = 16000
$ws.Range("O6").Value = 16000
$ws.Range("P6").Value = 16000
$ws.Range("S6").Value = 1067
$ws.Range("D7").Value = 44344
$ws.Range("L7").Value = "Segunda"
$ws.Range("M7").Value = 120
$ws.Range("N7").Value = 13000
$ws.Range("O7").Value = 13500
$ws.Range("P7").Value = 13250
$ws.Range("S7").Value = 883
$ws.Range("D8").Value = 44285
$ws.Range("M8").Value = 160
$ws.Range("N8").Value = 15000
$ws.Range("O8").Value = 16000
$ws.Range("P8").Value = 15500
$ws.Range("S8").Value = 1033
$ws.Range("D9").Value = 44295
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 160
$ws.Range("D10").Value = 44327
$ws.Range("N10").Value = 17000
$ws.Range("O10").Value = 17000
$ws.Range("P10").Value = 17000
$ws.Range("S10").Value = 1133
$ws.Range("D11").Value = 44327
$ws.Range("D12").Value = 44306
$ws.Range("M12").Value = 100
$ws.Range("N12").Value = 17500
$ws.Range("O12").Value = 17500
$ws.Range("P12").Value = 17500
$ws.Range("S12").Value = 1167
$ws.Range("D13").Value = 44306
$ws.Range("M13").Value = 200
$ws.Range("O13").Value = 14500
$ws.Range("P13").Value = 14250
$ws.Range("S13").Value = 950
$ws.Range("D14").Value = 44301
$ws.Range("D15").Value = 44301
$ws.Range("M15").Value = 80
$ws.Range("D16").Value = 44292
$ws.Range("L16").Value = "Segunda"
$ws.Range("M16").Value = 160
$ws.Range("N16").Value = 14000
$ws.Range("O16").Value = 15000
$ws.Range("P16").Value = 14500
$ws.Range("S16").Value = 967
$ws.Range("D17").Value = 44302
$ws.Range("L17").Value = "Primera"
$ws.Range("N17").Value = 17500
$ws.Range("O17").Value = 17500
$ws.Range("P17").Value = 17500
$ws.Range("S17").Value = 1167
$ws.Range("D18").Value = 44302
$ws.Range("L18").Value = "Segunda"
$ws.Range("M18").Value = 200
$ws.Range("N18").Value = 14000
$ws.Range("O18").Value = 15000
$ws.Range("P18").Value = 14500
$ws.Range("S18").Value = 967
$ws.Range("L19").Value = "Primera"
$ws.Range("M19").Value = 100
$ws.Range("N19").Value = 17500
$ws.Range("O19").Value = 17500
$ws.Range("P19").Value = 17500
$ws.Range("S19").Value = 1167
$ws.Range("D20").Value = 44316
$ws.Range("L20").Value = "Segunda"
$ws.Range("M20").Value = 200
$ws.Range("N20").Value = 14000
$ws.Range("O20").Value = 14500
$ws.Range("P20").Value = 14250
$ws.Range("Q20").Value = "$/caja 15 kilos empedrada"
$ws.Range("S20").Value = 950
$ws.Range("T20").Value = 15
$ws.Range("D21").Value = 44336
$ws.Range("M21").Value = 60
$ws.Range("N21").Value = 17000
$ws.Range("O21").Value = 17000
$ws.Range("P21").Value = 17000
$ws.Range("Q21").Value = "$/caja 15 kilos empedrada"
$ws.Range("S21").Value = 1133
$ws.Range("T21").Value = 15
$ws.Range("D22").Value = 44336
$ws.Range("M22").Value = 120
$ws.Range("O22").Value = 14500
$ws.Range("P22").Value = 14250
$ws.Range("Q22").Value = "$/caja 15 kilos empedrada"
$ws.Range("S22").Value = 950
$ws.Range("T22").Value = 15
$ws.Range("D23").Value = 44299
$ws.Range("M23").Value = 60
$ws.Range("D24").Value = 44299
$ws.Range("M24").Value = 120
$ws.Range("D25").Value = 44330
$ws.Range("N25").Value = 17000
$ws.Range("O25").Value = 17000
$ws.Range("P25").Value = 17000
$ws.Range("S25").Value = 1133
$ws.Range("D26").Value = 44330
$ws.Range("N26").Value = 14000
$ws.Range("O26").Value = 14500
$ws.Range("P26").Value = 14250
$ws.Range("S26").Value = 950
$ws.Range("D27").Value = 44298
$ws.Range("L27").Value = "Segunda"
$ws.Range("M27").Value = 80
$ws.Range("N27").Value = 14000
$ws.Range("P27").Value = 14500
$ws.Range("S27").Value = 967
$ws.Range("D28").Value = 44305
$ws.Range("L28").Value = "Primera"
$ws.Range("M28").Value = 60
$ws.Range("N28").Value = 17500
$ws.Range("O28").Value = 17500
$ws.Range("P28").Value = 17500
$ws.Range("S28").Value = 1167
$ws.Range("D29").Value = 44305
$ws.Range("L29").Value = "Segunda"
$ws.Range("M29").Value = 120
$ws.Range("N29").Value = 14000
$ws.Range("O29").Value = 15000
$ws.Range("P29").Value = 14500
$ws.Range("S29").Value = 967
$ws.Range("D30").Value = 44313
$ws.Range("L30").Value = "Especial"
$ws.Range("M30").Value = 100
$ws.Range("N30").Value = 17500
$ws.Range("O30").Value = 17500
$ws.Range("P30").Value = 17500
$ws.Range("Q30").Value = "$/caja 14 kilos empedrada"
$ws.Range("S30").Value = 1250
$ws.Range("T30").Value = 14
$ws.Range("D31").Value = 44313
$ws.Range("M31").Value = 100
$ws.Range("N31").Value = 16000
$ws.Range("O31").Value = 16000
$ws.Range("P31").Value = 16000
$ws.Range("Q31").Value = "$/caja 14 kilos empedrada"
$ws.Range("S31").Value = 1143
$ws.Range("T31").Value = 14
$ws.Range("D32").Value = 44313
$ws.Range("M32").Value = 80
$ws.Range("N32").Value = 14000
$ws.Range("O32").Value = 14000
$ws.Range("P32").Value = 14000
$ws.Range("Q32").Value = "$/caja 14 kilos empedrada"
$ws.Range("S32").Value = 1000
$ws.Range("T32").Value = 14
$ws.Range("D33").Value = 44323
$ws.Range("D34").Value = 44323
$ws.Range("M34").Value = 100
$ws.Range("O34").Value = 14000
$ws.Range("P34").Value = 14000
$ws.Range("S34").Value = 933
$ws.Range("D35").Value = 44348
$ws.Range("M35").Value = 100
$ws.Range("N35").Value = 15000
$ws.Range("O35").Value = 15000
$ws.Range("P35").Value = 15000
$ws.Range("S35").Value = 1000
$ws.Range("D36").Value = 44348
$ws.Range("M36").Value = 200
$ws.Range("N36").Value = 13000
$ws.Range("O36").Value = 13500
$ws.Range("P36").Value = 13250
$ws.Range("S36").Value = 883
$ws.Range("D37").Value = 44293
$ws.Range("N37").Value = 14000
$ws.Range("O37").Value = 15000
$ws.Range("P37").Value = 14500
$ws.Range("S37").Value = 967
$ws.Range("D38").Value = 44334
$ws.Range("L38").Value = "Primera"
$ws.Range("M38").Value = 200
$ws.Range("O38").Value = 17000
$ws.Range("P38").Value = 15500
$ws.Range("S38").Value = 1033
$ws.Range("D39").Value = 44334
$ws.Range("M39").Value = 100
$ws.Range("N39").Value = 14500
$ws.Range("O39").Value = 14500
$ws.Range("P39").Value = 14500
$ws.Range("S39").Value = 967
